$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the commessa value in row 2
$ws.Range("A2").Value = 254339

# Remove row 3 entirely (previously commessa 252980 with the same "motivo" text)
$ws.Rows("3:3").Delete()
